$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data set (rows 2-5), written column-by-column, top-to-bottom within
# each column, to match the order the shared-string table was rebuilt in.

# Column B - 종목명 (company name)
$ws.Cells.Item(2, 2).Value = "D-Wave Quantum Inc."
$ws.Cells.Item(3, 2).Value = "International Business Machines"
$ws.Cells.Item(4, 2).Value = "IonQ, Inc."
$ws.Cells.Item(5, 2).Value = "Rigetti Computing, Inc."

# Column C - 티커 (ticker)
$ws.Cells.Item(2, 3).Value = "QBTS"
$ws.Cells.Item(3, 3).Value = "IBM"
$ws.Cells.Item(4, 3).Value = "IONQ"
$ws.Cells.Item(5, 3).Value = "RGTI"

# Column D - 종가 (close price)
$ws.Cells.Item(2, 4).Value = 27.11
$ws.Cells.Item(3, 4).Value = 308.82
$ws.Cells.Item(4, 4).Value = 52.06
$ws.Cells.Item(5, 4).Value = 28.07

# Column E - RSI
$ws.Cells.Item(2, 5).Value = 59.6
$ws.Cells.Item(3, 5).Value = 52.9
$ws.Cells.Item(4, 5).Value = 57.7
$ws.Cells.Item(5, 5).Value = 56.4

# Column F - 5일수익률
$ws.Cells.Item(2, 6).Value = 19.59
$ws.Cells.Item(3, 6).Value = 0.08
$ws.Cells.Item(4, 6).Value = 5.59
$ws.Cells.Item(5, 6).Value = 9.779999999999999

# Column G - 점수(룰)
$ws.Cells.Item(2, 7).Value = 60
$ws.Cells.Item(3, 7).Value = 60
$ws.Cells.Item(4, 7).Value = 60
$ws.Cells.Item(5, 7).Value = 50

# Column H - 3일상승확률(%)
$ws.Cells.Item(2, 8).Value = 70
$ws.Cells.Item(3, 8).Value = 60
$ws.Cells.Item(4, 8).Value = 56
$ws.Cells.Item(5, 8).Value = 60

# Column I - 5일상승확률(%)
$ws.Cells.Item(2, 9).Value = 70
$ws.Cells.Item(3, 9).Value = 66
$ws.Cells.Item(4, 9).Value = 60
$ws.Cells.Item(5, 9).Value = 66

# Column J - 10일상승확률(%)
$ws.Cells.Item(2, 10).Value = 83
$ws.Cells.Item(3, 10).Value = 63
$ws.Cells.Item(4, 10).Value = 70
$ws.Cells.Item(5, 10).Value = 83

# Column K - 최종점수
$ws.Cells.Item(2, 11).Value = 61.2
$ws.Cells.Item(3, 11).Value = 59.6
$ws.Cells.Item(4, 11).Value = 57.2
$ws.Cells.Item(5, 11).Value = 56.6

# Column L (예측방식) is unchanged by this edit - left as-is.

# Column M - 판단 (judgment)
$ws.Cells.Item(2, 13).Value = "📈 매수 관찰 구간입니다."
$ws.Cells.Item(3, 13).Value = "⛔ 관망하십시오."
$ws.Cells.Item(4, 13).Value = "⛔ 관망하십시오."
$ws.Cells.Item(5, 13).Value = "⛔ 관망하십시오."

# Column N - MACRO_SCORE
$ws.Cells.Item(2, 14).Value = 50.60178744571824
$ws.Cells.Item(3, 14).Value = 50.60178744571824
$ws.Cells.Item(4, 14).Value = 50.60178744571824
$ws.Cells.Item(5, 14).Value = 50.60178744571824

# Column O - MACRO_SIGNAL
$ws.Cells.Item(2, 15).Value = "⚪ 중립 구간"
$ws.Cells.Item(3, 15).Value = "⚪ 중립 구간"
$ws.Cells.Item(4, 15).Value = "⚪ 중립 구간"
$ws.Cells.Item(5, 15).Value = "⚪ 중립 구간"
